# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the zh-cn and
# de-de sheets to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 07:19:09"
$wsZhCn.Range("H2").Value = "2016-03-24 07:19:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 07:19:14"
$wsDeDe.Range("H2").Value = "2016-03-24 07:19:40"
